$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 3 row update
$ws.Range("D9").Value = "Editor: Camera, Reload, Terxture, Entities"

# Week 2 row updates
$ws.Range("E6").Value = "finalize collision events and grounded"
$ws.Range("F6").ClearContents()
$ws.Range("E7").Value = "Zilch/JSON"
$ws.Range("F7").ClearContents()

# Update the active selection to match the author's final cursor position
$null = $ws.Range("F6").Select()
